# 自动更新价格数据: insert the newest day's prices at the top of the
# table (row 2), pushing all existing date rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 2..N down by inserting a fresh row above the current row 2.
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the header row by default;
# strip that back to the plain (unstyled) look used by the rest of the
# data rows.
$ws.Rows.Item(2).ClearFormats()

# Column A stores dates as plain text, not Excel date serials, so force
# a text number format before writing the value.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2026-01-02"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
